$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume-change (E) columns for each row.
# D-column values are forced to text (NumberFormat "@") before assignment
# and the cell style is reset to "Normal" afterwards so that numeric-looking
# strings (e.g. "584.14") are not silently converted into Excel numbers and
# no residual custom number format/style is left behind on the cell.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "70.722.75"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -2.29%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.628.78"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("E4").Value = "  +0.01%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "584.14"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.98%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "176.11"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -3.47%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "3.622.18"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.624"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +2.71%  "

$ws.Range("E10").Value = "  -4.99%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "6.85"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +17.59%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.611"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.52%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "48.55"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.69%  "

$ws.Range("E14").Value = "  -2.06%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "4.221.29"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.30%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "674.20"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -4.30%  "

$ws.Range("E17").Value = "  +0.29%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "3.633.52"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.32%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "70.788.36"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -2.24%  "

$ws.Range("E20").Value = "  -0.45%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "17.81"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -4.12%  "

$ws.Range("E22").Value = "  -1.79%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.946"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.25%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "17.23"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -3.55%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "100.11"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -4.28%  "

$ws.Range("E26").Value = "  -2.66%  "

$ws.Range("E27").Value = "  -2.89%  "

$ws.Range("E28").Value = "  -0.02%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.82"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -2.00%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "34.60"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -2.33%  "

$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("E32").Value = "  -5.66%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "7.57"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +1.52%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.38"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -6.75%  "

$ws.Range("E35").Value = "  -4.59%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "576.65"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -2.99%  "

$ws.Range("E37").Value = "  -2.03%  "

$ws.Range("E38").Value = "  -0.84%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "58.50"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -2.28%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.08%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.0453"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.57%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "3.544.65"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.39%  "

$ws.Range("E43").Value = "  -1.68%  "

$ws.Range("E44").Value = "  -3.62%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "34.43"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -4.11%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0733"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -6.03%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.98"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +4.88%  "

$ws.Range("E48").Value = "  -4.44%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.135"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +2.06%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "137.54"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +2.76%  "

$ws.Range("E51").Value = "  -2.55%  "
